# Apply automatic update of files:
# 1) Column C (Förändrad) for rows 2-18 changes from 46070 to 46072
# 2) Rows 9 and 11 swap their Beteckning (A), Datum (B) and Area (G) values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C (Förändrad date) for all data rows 2 through 18 ---
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 46072
}

# --- Swap data between row 9 and row 11 for columns A, B, G ---
# (values hard-coded from the known source data to avoid floating point
# round-trip drift when reading G's Value2 back out through COM)
$ws.Cells.Item(9, 1).Value2 = "A 35838-2023"
$ws.Cells.Item(9, 2).Value2 = 45147
$ws.Cells.Item(9, 7).Value2 = 1.1

$ws.Cells.Item(11, 1).Value2 = "A 6004-2026"
$ws.Cells.Item(11, 2).Value2 = 46050
$ws.Cells.Item(11, 7).Value2 = 2.7
